# ------------------------------------------------------------------
# PlayerPerformance_3125.xlsx edit:
#  - insert a new "Player Info" sheet before "ODI Batting"
#  - append a new "ODI Batting Extra" sheet after "ODI Bowling"
#  - on "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, replace
#    the scorecard URLs with bare match codes, drop the blank
#    INNING_NUMBER cells
#  - on "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE, replace
#    the scorecard URLs with bare match codes
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the value to be stored as text even when it looks numeric
    # (Excel would otherwise auto-coerce "2158" etc. into a number).
    if ($value -match '^-?[0-9]+(\.[0-9]+)?%?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
    $range.Borders.LineStyle = 1
}

# ------------------------------------------------------------------
# 1. Sheet restructuring
# ------------------------------------------------------------------

$wsBatting = $wb.Worksheets.Item("ODI Batting")

$wsInfo = $wb.Worksheets.Add($wsBatting)
$wsInfo.Name = "Player Info"

$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsBowling)
$wsExtra.Name = "ODI Batting Extra"

# ------------------------------------------------------------------
# 2. "Player Info" content
# ------------------------------------------------------------------

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $infoHeaders.Length; $i++) {
    $wsInfo.Cells.Item(1, $i + 1).Value = $infoHeaders[$i]
}
Set-HeaderStyle($wsInfo.Range("A1:D1"))

Set-TextValue $wsInfo.Cells.Item(2, 1) "3125"
$wsInfo.Cells.Item(2, 2).Value = "Fidel Henderson Edwards"
$wsInfo.Cells.Item(2, 3).Value = "Right Handed"
$wsInfo.Cells.Item(2, 4).Value = "Right Arm Fast"

Write-Host "Player Info sheet populated"

# ------------------------------------------------------------------
# 3. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, url -> bare code,
#    drop blank INNING_NUMBER cells
# ------------------------------------------------------------------

$wsBat = $wb.Worksheets.Item("ODI Batting")
$wsBat.Range("D1").Value = "MATCH_CODE"

$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $linkCell = $wsBat.Cells.Item($r, 4)
    $link = $linkCell.Value()
    if ($link -match 'MatchCode=(\d+)') {
        Set-TextValue $linkCell $matches[1]
    }

    $inningCell = $wsBat.Cells.Item($r, 2)
    $inning = $inningCell.Value()
    if ($inning -eq "") {
        $inningCell.ClearContents()
    }
}

Write-Host "ODI Batting sheet updated"

# ------------------------------------------------------------------
# 4. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, url -> bare code
# ------------------------------------------------------------------

$wsBowl = $wb.Worksheets.Item("ODI Bowling")
$wsBowl.Range("B1").Value = "MATCH_CODE"

$lastRowBowl = 50
for ($r = 2; $r -le $lastRowBowl; $r++) {
    $linkCell = $wsBowl.Cells.Item($r, 2)
    $link = $linkCell.Value()
    if ($link -match 'MatchCode=(\d+)') {
        Set-TextValue $linkCell $matches[1]
    }
}

Write-Host "ODI Bowling sheet updated"

# ------------------------------------------------------------------
# 5. "ODI Batting Extra" content
# ------------------------------------------------------------------

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $extraHeaders.Length; $i++) {
    $wsExtra.Cells.Item(1, $i + 1).Value = $extraHeaders[$i]
}
Set-HeaderStyle($wsExtra.Range("A1:F1"))

# columns: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraData = @(
    @("2776", "10", "1", "0", "4.73%", "NO"),
    @("2778", "10", "", "", "", "NO"),
    @("2780", "11", "2", "0", "4.56%", "NO"),
    @("2809", "10", "", "", "", "NO"),
    @("2811", "", "", "", "", "NO"),
    @("2825", "10", "1", "0", "3.14%", "NO"),
    @("2836", "", "", "", "", "NO"),
    @("2846", "", "", "", "", "NO"),
    @("2848", "", "", "", "", "NO"),
    @("2908", "", "", "", "", "NO"),
    @("2909", "11", "1", "0", "10.16%", "NO"),
    @("2910", "", "", "", "", "NO"),
    @("2911", "10", "0", "0", "0.34%", "NO"),
    @("2914", "10", "", "", "", "NO"),
    @("2948", "", "", "", "", "NO"),
    @("2949", "10", "0", "0", "", "NO"),
    @("2950", "", "", "", "", ""),
    @("2951", "", "", "", "", ""),
    @("2955", "", "", "", "", ""),
    @("2974", "", "", "", "", "")
)

for ($i = 0; $i -lt $extraData.Length; $i++) {
    $r = $i + 2
    $vals = $extraData[$i]

    Set-TextValue $wsExtra.Cells.Item($r, 1) $vals[0]

    if ($vals[1] -ne "") {
        $wsExtra.Cells.Item($r, 2).Value = [int]$vals[1]
    }

    if ($vals[2] -ne "") { Set-TextValue $wsExtra.Cells.Item($r, 3) $vals[2] }
    if ($vals[3] -ne "") { Set-TextValue $wsExtra.Cells.Item($r, 4) $vals[3] }
    if ($vals[4] -ne "") { Set-TextValue $wsExtra.Cells.Item($r, 5) $vals[4] }
    if ($vals[5] -ne "") { $wsExtra.Cells.Item($r, 6).Value = $vals[5] }
}

Write-Host "ODI Batting Extra sheet populated"

# ------------------------------------------------------------------
# 6. Restore the original active sheet
# ------------------------------------------------------------------

$wb.Worksheets.Item(1).Activate()
